$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new document row value for E112 (row 13)
$ws.Range("B13").Value = "Informe de Revision Tecnica Formal (RTF)"

# Move "Plan de Pruebas" from B15 (E114) up to B14 (E113)
$ws.Range("B14").Value = $ws.Range("B15").Value2
$ws.Range("B15").ClearContents()

# Update the active selection to C17, as recorded in the workbook view
$ws.Range("C17").Select()
